$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(15, 6).Value = 1348
$ws.Cells.Item(21, 6).Value = 24
$ws.Cells.Item(24, 6).Value = 2321
$ws.Cells.Item(29, 6).Value = 7
$ws.Cells.Item(34, 6).Value = 937

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 66

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(10, 6).Value = 3000
$ws.Cells.Item(11, 6).Value = 515
$ws.Cells.Item(14, 6).Value = 262

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(9, 6).Value = 3000
$ws.Cells.Item(10, 6).Value = 515
$ws.Cells.Item(30, 6).Value = 24
$ws.Cells.Item(35, 6).Value = 2321
$ws.Cells.Item(39, 6).Value = 262
$ws.Cells.Item(41, 6).Value = 7
$ws.Cells.Item(50, 6).Value = 937
